$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update trigger-condition wording ("，触发" inserted) for card effect descriptions in D2:D14.
# The underlying shared-string table is rebuilt by the engine in sheet scan order as a side
# effect of rewriting every cell top-to-bottom, which reproduces the target si ordering.
$ws.Range("D2").Value = "回合结束时：将主牌堆顶1张牌送墓。<br>`n开战时，触发：用墓地第1张牌替换本牌。"
$ws.Range("D3").Value = "开战时，触发：本牌所在槽位和对位槽位的所有怪物牌点数变为1。"
$ws.Range("D4").Value = "回合结束时，触发：如果本牌所在槽位和对位槽位的怪物牌合计数量不小于本牌点数，则将那些怪物牌全部消灭，然后消灭本牌。"
$ws.Range("D5").Value = "回合结束时，触发：横置本牌，然后本牌所在槽位和对位槽位中所有其他牌点数减1。"
$ws.Range("D6").Value = "回合结束时，触发：横置本牌所在槽位和对位槽位中所有牌。"
$ws.Range("D7").Value = "本牌所在行列的槽位新增怪物牌时，触发：那张怪物牌点数减2，然后本牌点数减1。"
$ws.Range("D8").Value = "死亡时，触发：如果本牌点数大于1，则将本牌移动到相邻槽位中而不是送墓，然后使本牌和那个槽位中所有其他牌点数减1。<br>"
$ws.Range("D9").Value = "开战时，触发：如果本牌所在槽位和对位槽位的怪物牌合计数量大于1，则将那些怪物牌洗回主牌堆，然后消灭本牌。"
$ws.Range("D10").Value = "点数降低时，触发：本牌所在槽位和相邻槽位的所有牌点数减1，然后消灭本牌。"
$ws.Range("D11").Value = "有牌进入本牌所在槽位时，触发：本牌所在槽位和对位槽位中所有牌点数减1，然后消灭本牌。"
$ws.Range("D12").Value = "死亡时，触发：选本牌所在行或列，其中的所有牌点数减1。"
$ws.Range("D13").Value = "回合结束时，触发：移动到1个相邻槽位，然后点数减1，本牌点数因此降至0时，消灭本牌所在槽位或对位槽位中的1张其他牌。"
$ws.Range("D14").Value = "回合结束时，触发：将本牌移动到对位槽位，同槽位中有怪物牌时，选其中1张一起移动到对位槽位。"

# Move the active selection to D16, matching the saved cursor position in the edited workbook.
$ws.Range("D16").Select()
